$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("J25").Value = "10-Susan"
$ws.Range("J26").Value = "20-David"
$ws.Range("J27").Value = "21-Daivd"
$ws.Range("J28").Value = "22b-Riaz"
$ws.Range("J29").Value = "22c-Riaz"
$ws.Range("J30").Value = "25-Dann"
$ws.Range("J31").Value = "26-Dann"
$ws.Range("J32").Value = "27-Dann"
$ws.Range("J33").Value = "28-Dann"
$ws.Range("J34").Value = "29-Dann"
$ws.Range("J35").Value = "30-Dann"
$ws.Range("J36").Value = "31-Dann"
$ws.Range("J37").Value = "32-Dann"

$ws.Range("N36").Select()
